# Applies the diff: re-maps/corrects species records in rows 10-24,
# adds two new observation rows (25-26), and extends the used range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 10-24 (records re-ordered / corrected) ---
# Row 10
$ws.Range("A10").Value = 112163866
$ws.Range("B10").Value = 90800
$ws.Range("E10").Value = 4364
$ws.Range("F10").Value = '''Dropptaggsvamp'
$ws.Range("G10").Value = '''Hydnellum ferrugineum'
$ws.Range("H10").Value = '''(Fr.:Fr.) P. Karst.'
$ws.Range("L10").Value = $null
$ws.Range("AI10").Value = '''Gammal barrblandskog, renlavsmarker'

# Row 11
$ws.Range("A11").Value = 112163367
$ws.Range("B11").Value = 90800
$ws.Range("P11").Value = '''Källdalen SV, Vstm'
$ws.Range("Q11").Value = 558083
$ws.Range("R11").Value = 6628611
$ws.Range("AI11").Value = '''Gammal barrblandskog'

# Row 12
$ws.Range("A12").Value = 112163316
$ws.Range("B12").Value = 90800
$ws.Range("D12").Value = '''LC'
$ws.Range("E12").Value = 4364
$ws.Range("F12").Value = '''Dropptaggsvamp'
$ws.Range("G12").Value = '''Hydnellum ferrugineum'
$ws.Range("H12").Value = '''(Fr.:Fr.) P. Karst.'
$ws.Range("I12").Value = '''21'
$ws.Range("P12").Value = '''Källdalen SV, Vstm'
$ws.Range("Q12").Value = 558054
$ws.Range("R12").Value = 6628665
$ws.Range("AC12").Value = '''Tätt bestånd.'
$ws.Range("AI12").Value = '''Gammal barrblandskog'

# Row 13
$ws.Range("A13").Value = 112164196
$ws.Range("B13").Value = 90823
$ws.Range("D13").Value = '''NT'
$ws.Range("E13").Value = 5966
$ws.Range("F13").Value = '''Motaggsvamp'
$ws.Range("G13").Value = '''Sarcodon squamosus'
$ws.Range("H13").Value = '''(Schaeff.) Quél.'
$ws.Range("I13").Value = '''1'
$ws.Range("J13").Value = '''fruktkroppar'
$ws.Range("P13").Value = '''Månses hål, Vstm'
$ws.Range("Q13").Value = 558007
$ws.Range("R13").Value = 6628378
$ws.Range("AI13").Value = '''Gammal barrblandskog, renlavsmarker'

# Row 14
$ws.Range("A14").Value = 112163912
$ws.Range("B14").Value = 90794
$ws.Range("E14").Value = 4362
$ws.Range("F14").Value = '''Blå taggsvamp'
$ws.Range("G14").Value = '''Hydnellum caeruleum'
$ws.Range("H14").Value = '''(Hornem.) P.Karst.'
$ws.Range("I14").Value = $null
$ws.Range("J14").Value = $null
$ws.Range("Q14").Value = 558020
$ws.Range("R14").Value = 6628293
$ws.Range("AC14").Value = $null
$ws.Range("AI14").Value = '''Gammal barrblandskog, renlavsmarker'
$ws.Range("AO14").Value = '''Gran'

# Row 15
$ws.Range("A15").Value = 112164100
$ws.Range("B15").Value = 95693
$ws.Range("E15").Value = 221941
$ws.Range("F15").Value = '''Plattlummer'
$ws.Range("G15").Value = '''Lycopodium complanatum'
$ws.Range("H15").Value = '''L.'
$ws.Range("L15").Value = $null
$ws.Range("Q15").Value = 558020
$ws.Range("R15").Value = 6628293
$ws.Range("AI15").Value = '''Gammal barrblandskog'

# Row 16
$ws.Range("A16").Value = 112164162
$ws.Range("B16").Value = 90800
$ws.Range("I16").Value = $null
$ws.Range("J16").Value = $null
$ws.Range("P16").Value = '''Månses hål, Vstm'
$ws.Range("Q16").Value = 558022
$ws.Range("R16").Value = 6628310
$ws.Range("AC16").Value = $null
$ws.Range("AI16").Value = '''Gammal barrblandskog, renlavsmarker'

# Row 17
$ws.Range("A17").Value = 112163592
$ws.Range("B17").Value = 90785
$ws.Range("E17").Value = 1968
$ws.Range("F17").Value = '''Grantaggsvamp'
$ws.Range("G17").Value = '''Bankera violascens'
$ws.Range("H17").Value = '''(Alb. & Schwein. : Fr.) Pouzar'
$ws.Range("I17").Value = '''2'
$ws.Range("J17").Value = '''fruktkroppar'
$ws.Range("Q17").Value = 558062
$ws.Range("R17").Value = 6628273
$ws.Range("AC17").Value = '''i bestånd med gamla granar.'
$ws.Range("AI17").Value = '''Gammal barrblandskog'
$ws.Range("AO17").Value = $null

# Row 18
$ws.Range("A18").Value = 112204239
$ws.Range("B18").Value = 89539
$ws.Range("P18").Value = '''Blåbärsberget SV, Vstm'
$ws.Range("Q18").Value = 558147
$ws.Range("R18").Value = 6627943
$ws.Range("AI18").Value = '''Blandskog'
$ws.Range("AO18").Value = '''Gran'

# Row 19
$ws.Range("A19").Value = 112204223
$ws.Range("B19").Value = 89503
$ws.Range("D19").Value = '''LC'
$ws.Range("E19").Value = 5447
$ws.Range("F19").Value = '''Vedticka'
$ws.Range("G19").Value = '''Fuscoporia viticola'
$ws.Range("H19").Value = '''(Schwein.) Murrill'
$ws.Range("P19").Value = '''Månses hål S, Vstm'
$ws.Range("Q19").Value = 558126
$ws.Range("R19").Value = 6627991
$ws.Range("AI19").Value = '''Barrblandskog'

# Row 20
$ws.Range("A20").Value = 112203759
$ws.Range("B20").Value = 89539
$ws.Range("D20").Value = '''NT'
$ws.Range("E20").Value = 1202
$ws.Range("F20").Value = '''Ullticka'
$ws.Range("G20").Value = '''Phellinidium ferrugineofuscum'
$ws.Range("H20").Value = '''(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q20").Value = 557997
$ws.Range("R20").Value = 6628183
$ws.Range("AO20").Value = '''Granlåga'

# Row 21
$ws.Range("A21").Value = 112203716
$ws.Range("B21").Value = 90800
$ws.Range("E21").Value = 4364
$ws.Range("F21").Value = '''Dropptaggsvamp'
$ws.Range("G21").Value = '''Hydnellum ferrugineum'
$ws.Range("H21").Value = '''(Fr.:Fr.) P. Karst.'
$ws.Range("L21").Value = $null
$ws.Range("M21").Value = $null
$ws.Range("Q21").Value = 558021
$ws.Range("R21").Value = 6628143
$ws.Range("S21").Value = 5
$ws.Range("AI21").Value = '''Barrblandskog, hedartad'
$ws.Range("AO21").Value = $null

# Row 22
$ws.Range("B22").Value = 95679

# Row 23
$ws.Range("A23").Value = 112203601
$ws.Range("B23").Value = 8377
$ws.Range("E23").Value = 106545
$ws.Range("F23").Value = '''Mindre märgborre'
$ws.Range("G23").Value = '''Tomicus minor'
$ws.Range("H23").Value = '''(Hartig, 1834)'
$ws.Range("L23").Value = $null
$ws.Range("M23").Value = '''äldre gnagspår'
$ws.Range("Q23").Value = 558087
$ws.Range("R23").Value = 6627982
$ws.Range("S23").Value = 25
$ws.Range("AI23").Value = '''Barrblandskog'
$ws.Range("AO23").Value = '''Tall'

# Row 24
$ws.Range("B24").Value = 89539

# --- Add two new rows (25-26) ---
# Row 25
$ws.Range("A25").Value = 112457954
$ws.Range("B25").Value = 89503
$ws.Range("C25").Value = '''Ovaliderad'
$ws.Range("D25").Value = '''LC'
$ws.Range("E25").Value = 5447
$ws.Range("F25").Value = '''Vedticka'
$ws.Range("G25").Value = '''Fuscoporia viticola'
$ws.Range("H25").Value = '''(Schwein.) Murrill'
$ws.Range("P25").Value = '''Månses hål, Vstm'
$ws.Range("Q25").Value = 558132
$ws.Range("R25").Value = 6628143
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = '''Västmanland'
$ws.Range("U25").Value = '''Surahammar'
$ws.Range("V25").Value = '''Västmanland'
$ws.Range("W25").Value = '''Ramnäs'
$ws.Range("Y25").Value = '''2023-05-02'
$ws.Range("AA25").Value = '''2023-05-02'
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AG25").Value = $false
$ws.Range("AI25").Value = '''Barrblandskog'
$ws.Range("AO25").Value = '''Gran'
$ws.Range("AW25").Value = '''Tom Sävström'
$ws.Range("AX25").Value = '''Tom Sävström'

# Row 26
$ws.Range("A26").Value = 112457599
$ws.Range("B26").Value = 96720
$ws.Range("C26").Value = '''Ovaliderad'
$ws.Range("D26").Value = '''VU'
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = '''Knärot'
$ws.Range("G26").Value = '''Goodyera repens'
$ws.Range("H26").Value = '''(L.) R. Br.'
$ws.Range("I26").Value = '''4'
$ws.Range("J26").Value = '''plantor/tuvor'
$ws.Range("K26").Value = '''fullt utvecklade blad'
$ws.Range("P26").Value = '''Månses hål, Vstm'
$ws.Range("Q26").Value = 558038
$ws.Range("R26").Value = 6628211
$ws.Range("S26").Value = 10
$ws.Range("T26").Value = '''Västmanland'
$ws.Range("U26").Value = '''Surahammar'
$ws.Range("V26").Value = '''Västmanland'
$ws.Range("W26").Value = '''Ramnäs'
$ws.Range("Y26").Value = '''2023-05-02'
$ws.Range("AA26").Value = '''2023-05-02'
$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AG26").Value = $false
$ws.Range("AI26").Value = '''Mossig gammal barrblandskog'
$ws.Range("AW26").Value = '''Tom Sävström'
$ws.Range("AX26").Value = '''Tom Sävström'
